$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells (F1:H1)
$ws.Range("F1").Value = "KNN_Outliers_MAD"
$ws.Range("G1").Value = "SVM_Outliers_MAD"
$ws.Range("H1").Value = "RF_Outliers_MAD"

# Match the formatting of the existing header cells (bold/centered/bordered style)
$ws.Range("E1").Copy()
$ws.Range("F1:H1").PasteSpecial(-4122)

# New boolean data cells (F2:H2)
$ws.Range("F2").Value = $false
$ws.Range("G2").Value = $false
$ws.Range("H2").Value = $false
